$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.102120194723294
$ws.Cells.Item(2, 4).Value = 1.099051961338214
$ws.Cells.Item(2, 5).Value = 1.114234865966422
$ws.Cells.Item(2, 6).Value = 1.117537558719009
$ws.Cells.Item(2, 9).Value = 1.076087562473893
$ws.Cells.Item(2, 10).Value = 1.106894590211531
$ws.Cells.Item(2, 11).Value = 1.101677628943686
$ws.Cells.Item(2, 12).Value = 1.116822843137496
$ws.Cells.Item(2, 13).Value = 1.120117483042885
$ws.Cells.Item(2, 14).Value = 1.108466507492216

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.103679922697071
$ws.Cells.Item(3, 4).Value = 1.10033974590183
$ws.Cells.Item(3, 5).Value = 1.115751141942976
$ws.Cells.Item(3, 6).Value = 1.119044790983635
$ws.Cells.Item(3, 9).Value = 1.076710851473056
$ws.Cells.Item(3, 10).Value = 1.10812088154234
$ws.Cells.Item(3, 11).Value = 1.10278600374313
$ws.Cells.Item(3, 12).Value = 1.118161923766851
$ws.Cells.Item(3, 13).Value = 1.121448126731981
$ws.Cells.Item(3, 14).Value = 1.109694540297372

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.10468655781842
$ws.Cells.Item(4, 4).Value = 1.101170204090461
$ws.Cells.Item(4, 5).Value = 1.116729913526745
$ws.Cells.Item(4, 6).Value = 1.120017638135905
$ws.Cells.Item(4, 9).Value = 1.077111342229048
$ws.Cells.Item(4, 10).Value = 1.10891128245145
$ws.Cells.Item(4, 11).Value = 1.103499812096899
$ws.Cells.Item(4, 12).Value = 1.119025516632918
$ws.Cells.Item(4, 13).Value = 1.122306185179208
$ws.Cells.Item(4, 14).Value = 1.110486063666433

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.105109132744103
$ws.Cells.Item(5, 4).Value = 1.101518661710576
$ws.Cells.Item(5, 5).Value = 1.117140834327761
$ws.Cells.Item(5, 6).Value = 1.120426050929197
$ws.Cells.Item(5, 9).Value = 1.077279039038809
$ws.Cells.Item(5, 10).Value = 1.10924283606839
$ws.Cells.Item(5, 11).Value = 1.103799095674669
$ws.Cells.Item(5, 12).Value = 1.119387890839413
$ws.Cells.Item(5, 13).Value = 1.122666214727977
$ws.Cells.Item(5, 14).Value = 1.110818088127538

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.105180049221648
$ws.Cells.Item(6, 4).Value = 1.101577130430754
$ws.Cells.Item(6, 5).Value = 1.11720979750407
$ws.Cells.Item(6, 6).Value = 1.120494591980725
$ws.Cells.Item(6, 9).Value = 1.077307156965545
$ws.Cells.Item(6, 10).Value = 1.109298462783838
$ws.Cells.Item(6, 11).Value = 1.103849299959338
$ws.Cells.Item(6, 12).Value = 1.119448695453857
$ws.Cells.Item(6, 13).Value = 1.122726624609362
$ws.Cells.Item(6, 14).Value = 1.110873793839303

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.104692206685987
$ws.Cells.Item(7, 4).Value = 1.10117486280772
$ws.Cells.Item(7, 5).Value = 1.116735406437727
$ws.Cells.Item(7, 6).Value = 1.120023097603054
$ws.Cells.Item(7, 9).Value = 1.077113585626122
$ws.Cells.Item(7, 10).Value = 1.108915715547764
$ws.Cells.Item(7, 11).Value = 1.10350381427707
$ws.Cells.Item(7, 12).Value = 1.119030361355681
$ws.Cells.Item(7, 13).Value = 1.122310998643214
$ws.Cells.Item(7, 14).Value = 1.110490503058253

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.102647859888752
$ws.Cells.Item(8, 4).Value = 1.099487764558405
$ws.Cells.Item(8, 5).Value = 1.114747793487073
$ws.Cells.Item(8, 6).Value = 1.118047444890255
$ws.Cells.Item(8, 9).Value = 1.076298792814556
$ws.Cells.Item(8, 10).Value = 1.107309666804293
$ws.Cells.Item(8, 11).Value = 1.102052915413728
$ws.Cells.Item(8, 12).Value = 1.117275994263278
$ws.Cells.Item(8, 13).Value = 1.120567798428859
$ws.Cells.Item(8, 14).Value = 1.108882173541349

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.099024944878801
$ws.Cells.Item(9, 4).Value = 1.09649283470093
$ws.Cells.Item(9, 5).Value = 1.111226823637951
$ws.Cells.Item(9, 6).Value = 1.114546995310073
$ws.Cells.Item(9, 9).Value = 1.074841193180937
$ws.Cells.Item(9, 10).Value = 1.104455510953423
$ws.Cells.Item(9, 11).Value = 1.099469932736697
$ws.Cells.Item(9, 12).Value = 1.114162067321344
$ws.Cells.Item(9, 13).Value = 1.117472973176619
$ws.Cells.Item(9, 14).Value = 1.106023964461795

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.096595171219803
$ws.Cells.Item(10, 4).Value = 1.094480808266609
$ws.Cells.Item(10, 5).Value = 1.108866391013653
$ws.Cells.Item(10, 6).Value = 1.112199870146776
$ws.Cells.Item(10, 9).Value = 1.073854443032921
$ws.Cells.Item(10, 10).Value = 1.102535958225246
$ws.Cells.Item(10, 11).Value = 1.097729698252792
$ws.Cells.Item(10, 12).Value = 1.1120703691023
$ws.Cells.Item(10, 13).Value = 1.115393622038903
$ws.Cells.Item(10, 14).Value = 1.104101685748568

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.095539456243773
$ws.Cells.Item(11, 4).Value = 1.093605793593573
$ws.Cells.Item(11, 5).Value = 1.107841039563908
$ws.Cells.Item(11, 6).Value = 1.111180192693464
$ws.Cells.Item(11, 9).Value = 1.073423531310431
$ws.Cells.Item(11, 10).Value = 1.101700659384212
$ws.Cells.Item(11, 11).Value = 1.096971705521693
$ws.Cells.Item(11, 12).Value = 1.111160768129803
$ws.Cells.Item(11, 13).Value = 1.11448927554915
$ws.Cells.Item(11, 14).Value = 1.10326520068737

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.09514676110674
$ws.Cells.Item(12, 4).Value = 1.093280192810235
$ws.Cells.Item(12, 5).Value = 1.107459674608827
$ws.Cells.Item(12, 6).Value = 1.110800922068431
$ws.Cells.Item(12, 9).Value = 1.073262917638531
$ws.Cells.Item(12, 10).Value = 1.101389761866981
$ws.Cells.Item(12, 11).Value = 1.096689472436916
$ws.Cells.Item(12, 12).Value = 1.110822306665033
$ws.Cells.Item(12, 13).Value = 1.114152752022571
$ws.Cells.Item(12, 14).Value = 1.102953861660003

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.095231020909518
$ws.Cells.Item(13, 4).Value = 1.093350061775877
$ws.Cells.Item(13, 5).Value = 1.10754150168857
$ws.Cells.Item(13, 6).Value = 1.110882300506475
$ws.Cells.Item(13, 9).Value = 1.073297394987653
$ws.Cells.Item(13, 10).Value = 1.101456479147769
$ws.Cells.Item(13, 11).Value = 1.096750043375709
$ws.Cells.Item(13, 12).Value = 1.110894934923966
$ws.Cells.Item(13, 13).Value = 1.114224965212242
$ws.Cells.Item(13, 14).Value = 1.103020673686983

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.095507007354772
$ws.Cells.Item(14, 4).Value = 1.093578891257335
$ws.Cells.Item(14, 5).Value = 1.107809526168571
$ws.Cells.Item(14, 6).Value = 1.111148852686238
$ws.Cells.Item(14, 9).Value = 1.073410266271705
$ws.Cells.Item(14, 10).Value = 1.101674973423352
$ws.Cells.Item(14, 11).Value = 1.096948389998641
$ws.Cells.Item(14, 12).Value = 1.111132803010338
$ws.Cells.Item(14, 13).Value = 1.114461470902029
$ws.Cells.Item(14, 14).Value = 1.103239478249499

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.095676977679567
$ws.Cells.Item(15, 4).Value = 1.093719803273135
$ws.Cells.Item(15, 5).Value = 1.107974597759288
$ws.Cells.Item(15, 6).Value = 1.1113130153894
$ws.Cells.Item(15, 9).Value = 1.073479736401698
$ws.Cells.Item(15, 10).Value = 1.101809511089736
$ws.Cells.Item(15, 11).Value = 1.097070507356642
$ws.Cells.Item(15, 12).Value = 1.111279282164054
$ws.Cells.Item(15, 13).Value = 1.11460710880802
$ws.Cells.Item(15, 14).Value = 1.103374206974801

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.096665158386669
$ws.Cells.Item(16, 4).Value = 1.094538799206305
$ws.Cells.Item(16, 5).Value = 1.108934370266072
$ws.Cells.Item(16, 6).Value = 1.112267470976812
$ws.Cells.Item(16, 9).Value = 1.073882963918271
$ws.Cells.Item(16, 10).Value = 1.102591306489713
$ws.Cells.Item(16, 11).Value = 1.097779908852546
$ws.Cells.Item(16, 12).Value = 1.112130653519183
$ws.Cells.Item(16, 13).Value = 1.11545355579665
$ws.Cells.Item(16, 14).Value = 1.104157112613919

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.097284042736935
$ws.Cells.Item(17, 4).Value = 1.095051509976554
$ws.Cells.Item(17, 5).Value = 1.109535526188412
$ws.Cells.Item(17, 6).Value = 1.112865268123239
$ws.Cells.Item(17, 9).Value = 1.074134918045741
$ws.Cells.Item(17, 10).Value = 1.103080595096823
$ws.Cells.Item(17, 11).Value = 1.098223696078308
$ws.Cells.Item(17, 12).Value = 1.112663648593778
$ws.Cells.Item(17, 13).Value = 1.115983437341798
$ws.Cells.Item(17, 14).Value = 1.104647096066973

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.097644680421539
$ws.Cells.Item(18, 4).Value = 1.095350200208234
$ws.Cells.Item(18, 5).Value = 1.109885855765497
$ws.Cells.Item(18, 6).Value = 1.113213630070369
$ws.Cells.Item(18, 9).Value = 1.074281527639606
$ws.Cells.Item(18, 10).Value = 1.103365592132995
$ws.Cells.Item(18, 11).Value = 1.098482120097344
$ws.Cells.Item(18, 12).Value = 1.112974161876557
$ws.Cells.Item(18, 13).Value = 1.116292125690361
$ws.Cells.Item(18, 14).Value = 1.104932497831626

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.097767590092253
$ws.Cells.Item(19, 4).Value = 1.095451984241721
$ws.Cells.Item(19, 5).Value = 1.110005256135418
$ws.Cells.Item(19, 6).Value = 1.113332358085122
$ws.Cells.Item(19, 9).Value = 1.074331458388904
$ws.Cells.Item(19, 10).Value = 1.103462701807272
$ws.Cells.Item(19, 11).Value = 1.098570163473807
$ws.Cells.Item(19, 12).Value = 1.113079975850924
$ws.Cells.Item(19, 13).Value = 1.116397315905818
$ws.Cells.Item(19, 14).Value = 1.10502974541278

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.097217678302862
$ws.Cells.Item(20, 4).Value = 1.094996538828622
$ws.Cells.Item(20, 5).Value = 1.109471060469463
$ws.Cells.Item(20, 6).Value = 1.112801163649607
$ws.Cells.Item(20, 9).Value = 1.074107922104456
$ws.Cells.Item(20, 10).Value = 1.103028140177732
$ws.Cells.Item(20, 11).Value = 1.098176126418628
$ws.Cells.Item(20, 12).Value = 1.11260650200922
$ws.Cells.Item(20, 13).Value = 1.115926625726492
$ws.Cells.Item(20, 14).Value = 1.10459456665588

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.095425751673032
$ws.Cells.Item(21, 4).Value = 1.093511522819197
$ws.Cells.Item(21, 5).Value = 1.107730613648894
$ws.Cells.Item(21, 6).Value = 1.11107037408783
$ws.Cells.Item(21, 9).Value = 1.073377043835071
$ws.Cells.Item(21, 10).Value = 1.101610649797343
$ws.Cells.Item(21, 11).Value = 1.096890000743631
$ws.Cells.Item(21, 12).Value = 1.111062773313727
$ws.Cells.Item(21, 13).Value = 1.114391842778545
$ws.Cells.Item(21, 14).Value = 1.103175063276562

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.094295874929246
$ws.Cells.Item(22, 4).Value = 1.092574464941935
$ws.Cells.Item(22, 5).Value = 1.106633404586854
$ws.Cells.Item(22, 6).Value = 1.109979160301807
$ws.Cells.Item(22, 9).Value = 1.072914303785145
$ws.Cells.Item(22, 10).Value = 1.100715764683209
$ws.Cells.Item(22, 11).Value = 1.096077417946469
$ws.Cells.Item(22, 12).Value = 1.110088720363053
$ws.Cells.Item(22, 13).Value = 1.113423334467284
$ws.Cells.Item(22, 14).Value = 1.102278907322922

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.09489515392162
$ws.Cells.Item(23, 4).Value = 1.093071540258958
$ws.Cells.Item(23, 5).Value = 1.107215337031973
$ws.Cells.Item(23, 6).Value = 1.110557921797723
$ws.Cells.Item(23, 9).Value = 1.073159917346282
$ws.Cells.Item(23, 10).Value = 1.101190510272134
$ws.Cells.Item(23, 11).Value = 1.096508560997859
$ws.Cells.Item(23, 12).Value = 1.110605415255548
$ws.Cells.Item(23, 13).Value = 1.113937097617776
$ws.Cells.Item(23, 14).Value = 1.102754327105038

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.097247666606437
$ws.Cells.Item(24, 4).Value = 1.095021379052875
$ws.Cells.Item(24, 5).Value = 1.109500190724758
$ws.Cells.Item(24, 6).Value = 1.112830130699995
$ws.Cells.Item(24, 9).Value = 1.074120121492796
$ws.Cells.Item(24, 10).Value = 1.103051843523618
$ws.Cells.Item(24, 11).Value = 1.098197622427504
$ws.Cells.Item(24, 12).Value = 1.112632325245712
$ws.Cells.Item(24, 13).Value = 1.115952297632188
$ws.Cells.Item(24, 14).Value = 1.104618303663236

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.099964058405922
$ws.Cells.Item(25, 4).Value = 1.097269768238825
$ws.Cells.Item(25, 5).Value = 1.112139342858918
$ws.Cells.Item(25, 6).Value = 1.115454275285395
$ws.Cells.Item(25, 9).Value = 1.075220639386138
$ws.Cells.Item(25, 10).Value = 1.105196293482221
$ws.Cells.Item(25, 11).Value = 1.100140870677676
$ws.Cells.Item(25, 12).Value = 1.114969822198621
$ws.Cells.Item(25, 13).Value = 1.118275858068329
$ws.Cells.Item(25, 14).Value = 1.106765798986753
